# Update countries & provincias Spain
#
# 1) Refresh the "last updated" timestamp in A1 (13:22 -> 13:52).
# 2) "Principado de Andorra" jumps up the (case-count sorted) table from
#    row 89 to row 85 (just after "Libano") with refreshed numbers; the
#    countries that used to sit at rows 85-88 (Banglades, Cuba, Republica
#    de Chipre, Afganistan) simply slide down one row each, keeping their
#    existing figures.
# 3) Same pattern for "Mozambique": it moves from row 168 up to row 165
#    (just after "Libia") with refreshed numbers, and Antigua y Barbuda,
#    Somalia, Guinea Ecuatorial slide down one row each.
# 4) A handful of other countries simply received refreshed case counts
#    in place (Noruega, Barein, Bosnia y Herzegovina, Senegal).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) timestamp ---------------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 12 de Abril de 2020 a las 13:52"

# --- helper: write a full data row (B..H) ---------------------------------
function Set-Row($row, $b, $c, $d, $e, $f, $g, $h) {
    $ws.Range("B$row").Value = $b
    $ws.Range("C$row").Value = $c
    $ws.Range("D$row").Value = $d
    $ws.Range("E$row").Value = $e
    $ws.Range("F$row").Value = $f
    $ws.Range("G$row").Value = $g
    $ws.Range("H$row").Value = $h
}

# --- 2) simple in-place refreshes -----------------------------------------
# Noruega (row 30)
Set-Row 30 6459 50 32 6304 67 4 123

# Barein (row 71)
Set-Row 71 1024 0 557 461 3 0 6

# Bosnia y Herzegovina (row 74)
Set-Row 74 986 40 193 755 4 1 38

# Senegal (row 108)
Set-Row 108 280 2 171 107 1 0 2

# --- 3) Andorra moves from row 89 to row 85 --------------------------------
# Slide Banglades/Cuba/Republica de Chipre/Afganistan (rows 85-88) down
# to rows 86-89, then drop the refreshed Andorra figures into row 85.
$ws.Range("A85:H88").Copy($ws.Range("A86:H89"))
$ws.Range("A85").Value = "Principado de Andorra"
Set-Row 85 622 21 71 523 17 2 28

# --- 4) Mozambique moves from row 168 to row 165 ---------------------------
# Slide Antigua y Barbuda/Somalia/Guinea Ecuatorial (rows 165-167) down
# to rows 166-168, then drop the refreshed Mozambique figures into row 165.
$ws.Range("A165:H167").Copy($ws.Range("A166:H168"))
$ws.Range("A165").Value = "Mozambique"
Set-Row 165 21 1 2 19 0 0 0
